$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0496
$ws.Range("E2").Value = -0.243
$ws.Range("F2").Value = 0.4
$ws.Range("G2").Value = 0.04271076341404553
$ws.Range("H2").Value = 0.04271076341404553
$ws.Range("I2").Value = 0.02096070308100716
$ws.Range("J2").Value = 0.01559766017011187
$ws.Range("K2").Value = 183.2
$ws.Range("L2").Value = 0.009210843858096694
$ws.Range("M2").Value = 50.5
$ws.Range("N2").Value = 0.008408958454749812
$ws.Range("O2").Value = 0.2756550218340612
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 50.5
$ws.Range("T2").Value = 1
$ws.Range("U2").Value = 2541.1
$ws.Range("V2").Value = 0.4231287986012822
$ws.Range("W2").Value = 0.02598802734984538
$ws.Range("X2").Value = 0.1049515084394596
$ws.Range("Y2").Value = -0.07896348108961419
$ws.Range("Z2").Value = 2.743809405564982
$ws.Range("AA2").Value = 0.04279700667955924
$ws.Range("AB2").Value = 0.07014030709155014
$ws.Range("AC2").Value = -0.0273433004119909
$ws.Range("AD2").Value = 3911.5
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 3911.5
$ws.Range("AG2").Value = 1370.4
$ws.Range("AH2").Value = 0.3944237168498538
$ws.Range("AI2").Value = 0.3477043424152185
$ws.Range("AJ2").Value = 0.1857942759527651
$ws.Range("AK2").Value = 0.157365302466584
$ws.Range("AL2").Value = 132.7
$ws.Range("AM2").Value = 132.7
$ws.Range("AN2").Value = 8.064948453608247
$ws.Range("AO2").Value = 3.141672946495855
$ws.Range("AP2").Value = 2.825567010309279
$ws.Range("AQ2").Value = 3.141672946495855

$ws.Range("D3").Value = 0.0496
$ws.Range("E3").Value = -0.243
$ws.Range("F3").Value = 0.4
$ws.Range("G3").Value = 0.04271076341404553
$ws.Range("H3").Value = 0.04271076341404553
$ws.Range("I3").Value = 0.02096070308100716
$ws.Range("J3").Value = 0.01559766017011187
$ws.Range("K3").Value = 183.2
$ws.Range("L3").Value = 0.009210843858096694
$ws.Range("M3").Value = 50.5
$ws.Range("N3").Value = 0.008408958454749812
$ws.Range("O3").Value = 0.2756550218340612
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 50.5
$ws.Range("T3").Value = 1
$ws.Range("U3").Value = 2541.1
$ws.Range("V3").Value = 0.4231287986012822
$ws.Range("W3").Value = 0.02598802734984538
$ws.Range("X3").Value = 0.1049515084394596
$ws.Range("Y3").Value = -0.07896348108961419
$ws.Range("Z3").Value = 2.743809405564982
$ws.Range("AA3").Value = 0.04279700667955924
$ws.Range("AB3").Value = 0.07014030709155014
$ws.Range("AC3").Value = -0.0273433004119909
$ws.Range("AD3").Value = 3911.5
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 3911.5
$ws.Range("AG3").Value = 1370.4
$ws.Range("AH3").Value = 0.3944237168498538
$ws.Range("AI3").Value = 0.3477043424152185
$ws.Range("AJ3").Value = 0.1857942759527651
$ws.Range("AK3").Value = 0.157365302466584
$ws.Range("AL3").Value = 132.7
$ws.Range("AM3").Value = 132.7
$ws.Range("AN3").Value = 8.064948453608247
$ws.Range("AO3").Value = 3.141672946495855
$ws.Range("AP3").Value = 2.825567010309279
$ws.Range("AQ3").Value = 3.141672946495855
